$wb = $excel.ActiveWorkbook

# --- Restricciones_del_lider ---
$ws = $wb.Worksheets.Item("Restricciones_del_lider")
$ws.Range("A2").Value = "1.9399999999999995 - x"
$ws.Range("B2").Value = "-2.9399999999999995"
$ws.Range("D2").Value = "0.32"
$ws.Range("A3").Value = "-1.9399999999999997 + x"
$ws.Range("B3").Value = "0.9399999999999997"
$ws.Range("D3").Value = "0.02"
$ws.Range("A4").Value = "35.63239999999999 + x - y - 9(x^2)"
$ws.Range("B4").Value = "-34.63239999999999"
$ws.Range("D4").Value = "0.44"

# --- Restricciones_del_follower ---
$ws = $wb.Worksheets.Item("Restricciones_del_follower")
$ws.Range("A2").Value = "-19.71359999999999 + (-0.5 + x)*(y^2)"
$ws.Range("B2").Value = "19.71359999999999"
$ws.Range("D2").Value = "0.36"
$ws.Range("F2").Value = "0"
$ws.Range("A3").Value = "-3.6999999999999993 + y"
$ws.Range("B3").Value = "2.6999999999999993"
$ws.Range("D3").Value = "0.43"
$ws.Range("E3").Value = "-9.0"
$ws.Range("F3").Value = "-0.1"
$ws.Range("A4").Value = "-5.699999999999999 - y"
$ws.Range("B4").Value = "-4.699999999999999"
$ws.Range("D4").Value = "0.43"
$ws.Range("E4").Value = "0"
$ws.Range("F4").Value = "0"

# --- Punto_modificado ---
$ws = $wb.Worksheets.Item("Punto_modificado")
$ws.Range("A2").Value = "1.9399999999999997"
$ws.Range("B2").Value = "3.6999999999999993"

# --- Vector_bf ---
$ws = $wb.Worksheets.Item("Vector_bf")
$ws.Range("A2").Value = "-4.836159999999998"

# --- Vector_BF ---
$ws = $wb.Worksheets.Item("Vector_BF")
$ws.Range("A2").Value = "11.3448"
$ws.Range("A3").Value = "9.44"
